# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a 🚀
#
# Updates the "Metadata" sheet (version/date/publisher/jurisdiction) and the
# "Elements" sheet (root Extension row's Short/Definition -> Problem Type row)
# to match the new IG build.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Version bump: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Build date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher was blank, now populated
$meta.Range("B9").Value = "Alvearie Team"

# First "Contact" row becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# The second (duplicate) "Contact" / "No display for ContactDetail" row is removed entirely
$meta.Rows.Item(11).Delete()

# Elements sheet: the root Extension row's Short/Definition get replaced with
# the ProblemType-specific text
$elements.Range("K2").Value = "Problem Type"
$elements.Range("L2").Value = "Problem type code"
